$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update the English (column C) translations that were re-worded.
# ---------------------------------------------------------------------------
$ws.Range("C2").Value = "Security Group is a distributed and stateless virtual firewall. It can be used to control network "
$ws.Range("C3").Value = "access for virtual machines or containers."
$ws.Range("C10").Value = "The number of selected security groups has reached the upper limit."
$ws.Range("C11").Value = "At least one security group is required."

# ---------------------------------------------------------------------------
# 2. Re-apply explicit font formatting on B2 (size 12 Arial) - this is what
#    produces the new dedicated cell style used only by that cell.
# ---------------------------------------------------------------------------
$ws.Range("B2").Font.Name = "Arial"
$ws.Range("B2").Font.Size = 12

# ---------------------------------------------------------------------------
# 3. Re-apply word-wrap to C10/C11 so they pick up their own (duplicate)
#    wrap-text style, matching the new cellXfs layout.
# ---------------------------------------------------------------------------
$ws.Range("C10:C11").WrapText = $true

# ---------------------------------------------------------------------------
# 4. Row height / column width adjustments.
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 30
$ws.Columns.Item(2).ColumnWidth = 34.21428571
$ws.Columns.Item(3).ColumnWidth = 36.78571428

# ---------------------------------------------------------------------------
# 5. Page setup (paper size / orientation) for printing.
# ---------------------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# 6. Selection / scroll position bookkeeping to match the saved view state.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("C15").Select()

# ---------------------------------------------------------------------------
# 7. Disable iterative calculation (workbook previously had iterate="1").
# ---------------------------------------------------------------------------
$excel.Iteration = $false
